# Se actualiza los datos de pagos
$wb = $excel.ActiveWorkbook

# Update the shared payment data values (numeroUsuario/placa/fecha columns,
# cell A2 / C2 / G2 (or E2) on each "Pago..." sheet) -- these three values
# share the same underlying strings across sheets 1-5.
foreach ($sheetName in @("PagoSinTarjetaAsociada","PagoAfiliadoDebitoAuto","PagoConValidacionHistorial","PagoAsociandoTarjeta")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A2").Value = "'72934725"
    $ws.Range("C2").Value = "'XFN-363"
    $ws.Range("G2").Value = "'06/11/2025"
}

# PagoConTarjetaAsociada has a narrower table (only columns A-E), so the
# "fecha" equivalent value lands in E2 instead of G2.
$ws4 = $wb.Worksheets.Item("PagoConTarjetaAsociada")
$ws4.Range("A2").Value = "'72934725"
$ws4.Range("C2").Value = "'XFN-363"
$ws4.Range("E2").Value = "'06/11/2025"

# Move the active/selected tab from "PagoSinTarjetaAsociada" to
# "PagoAsociandoTarjeta", and move the selection on the now-inactive first
# sheet from A2 to G2.
$wsOld = $wb.Worksheets.Item("PagoSinTarjetaAsociada")
$wsOld.Range("G2").Select() | Out-Null

$wsNew = $wb.Worksheets.Item("PagoAsociandoTarjeta")
$wsNew.Activate()
$wsNew.Range("G2").Select() | Out-Null
